# adding averages and more checks
$wb = $excel.ActiveWorkbook

$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1) Header / title font styling: the big title (row 1) drops its 14pt size
#    (reverting to the default size) and, together with the column-header
#    row (row 2), switches its bold font color to white so it shows up
#    against the dark-blue header fill.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsTraining, $wsExam)) {
    $ws.Range("A1").Font.Size = 11
    $ws.Range("A1").Font.Color = 16777215
}
$wsTraining.Range("A2:K2").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# 2) Training Dashboard: "PERIOD TO EXPIRE" (col H) shrinks by 8 (8 more
#    days elapsed) and "LAST UPDATE" (col I) moves from 08-Sep-2025 to
#    16-Sep-2025 for every data row (3-27).
# ---------------------------------------------------------------------------
$periodToExpire = @{
    3 = 426; 4 = 213; 5 = 310; 6 = 240; 7 = 408; 8 = 525; 9 = 219; 10 = 247;
    11 = 489; 12 = 400; 13 = 210; 14 = 255; 15 = 399; 16 = 258; 17 = 489;
    18 = 82; 19 = -103; 20 = -343; 21 = -41; 22 = -41; 23 = 155; 24 = 278;
    25 = 314; 26 = 314; 27 = 308
}

foreach ($row in $periodToExpire.Keys) {
    $wsTraining.Cells.Item($row, 8).Value = $periodToExpire[$row]
}

# Force the "LAST UPDATE" column to keep storing plain text (not an
# auto-converted date serial) while writing the new date string.
$wsTraining.Range("I3:I27").NumberFormat = "@"
for ($row = 3; $row -le 27; $row++) {
    $wsTraining.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 3) Exam Dashboard: COMMENTS column (E) is narrowed and every remark is
#    replaced with "date is valid".
# ---------------------------------------------------------------------------
$wsExam.Columns.Item(5).ColumnWidth = 14.14

for ($row = 3; $row -le 8; $row++) {
    $wsExam.Cells.Item($row, 5).Value = "date is valid"
}
